$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44176
$ws.Range("J2").Value = 700

$ws.Range("D3").Value = 44491
$ws.Range("J3").Value = 500

$ws.Range("D4").Value = 44482
$ws.Range("J4").Value = 72
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 2000
$ws.Range("P4").Value = 2000

$ws.Range("D5").Value = 44497
$ws.Range("J5").Value = 120

$ws.Range("D6").Value = 44537
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("J6").Value = 480

$ws.Range("D7").Value = 44519
$ws.Range("J7").Value = 600
$ws.Range("L7").Value = 1800
$ws.Range("M7").Value = 1700
$ws.Range("P7").Value = 1700

$ws.Range("D8").Value = 44475
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 2000
$ws.Range("P8").Value = 2000

$ws.Range("D9").Value = 44169
$ws.Range("H9").Value = "Verde"
$ws.Range("J9").Value = 600
$ws.Range("K9").Value = 1600
$ws.Range("L9").Value = 1600
$ws.Range("M9").Value = 1600
$ws.Range("P9").Value = 1600

$ws.Range("D10").Value = 44544
$ws.Range("H10").Value = "Sin especificar"

$ws.Range("D11").Value = 44539
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("J11").Value = 120

$ws.Range("D12").Value = 44474
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = 2000
$ws.Range("P12").Value = 2000

$ws.Range("D13").Value = 44516
$ws.Range("J13").Value = 360
$ws.Range("K13").Value = 1600
$ws.Range("L13").Value = 1600
$ws.Range("M13").Value = 1600
$ws.Range("P13").Value = 1600

$ws.Range("D14").Value = 44523
$ws.Range("J14").Value = 520
$ws.Range("K14").Value = 1800
$ws.Range("L14").Value = 1800
$ws.Range("M14").Value = 1800
$ws.Range("P14").Value = 1800

$ws.Range("D15").Value = 44488
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 1700
$ws.Range("L15").Value = 1800
$ws.Range("M15").Value = 1750
$ws.Range("P15").Value = 1750

$ws.Range("D16").Value = 44530
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = 1500
$ws.Range("P16").Value = 1500

$ws.Range("D17").Value = 44540
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("J17").Value = 600
$ws.Range("K17").Value = 1700
$ws.Range("L17").Value = 1700
$ws.Range("M17").Value = 1700
$ws.Range("P17").Value = 1700

$ws.Range("D18").Value = 44494
$ws.Range("K18").Value = 1700
$ws.Range("L18").Value = 1700
$ws.Range("M18").Value = 1700
$ws.Range("P18").Value = 1700

$ws.Range("D19").Value = 44518
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = 1600
$ws.Range("L19").Value = 1600
$ws.Range("M19").Value = 1600
$ws.Range("P19").Value = 1600

$ws.Range("D20").Value = 44509
$ws.Range("J20").Value = 550

$ws.Range("D21").Value = 44165
$ws.Range("H21").Value = "Verde"
$ws.Range("J21").Value = 300

$ws.Range("D22").Value = 44168
$ws.Range("H22").Value = "Verde"
$ws.Range("J22").Value = 200

$ws.Range("D23").Value = 44532
$ws.Range("J23").Value = 180
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 1500
$ws.Range("M23").Value = 1500
$ws.Range("P23").Value = 1500

$ws.Range("D24").Value = 44533
$ws.Range("J24").Value = 420
$ws.Range("K24").Value = 1700
$ws.Range("L24").Value = 1700
$ws.Range("M24").Value = 1700
$ws.Range("P24").Value = 1700

$ws.Range("D25").Value = 44473
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 2000
$ws.Range("P25").Value = 2000

$ws.Range("D26").Value = 44504
$ws.Range("J26").Value = 180
$ws.Range("K26").Value = 1600
$ws.Range("L26").Value = 1600
$ws.Range("M26").Value = 1600
$ws.Range("P26").Value = 1600

$ws.Range("D27").Value = 44525
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("J27").Value = 180

$ws.Range("D28").Value = 44487
$ws.Range("K28").Value = 1800
$ws.Range("L28").Value = 1800
$ws.Range("M28").Value = 1800
$ws.Range("P28").Value = 1800

$ws.Range("D29").Value = 44179
$ws.Range("H29").Value = "Verde"
$ws.Range("J29").Value = 200

$ws.Range("D30").Value = 44503
$ws.Range("J30").Value = 72

$ws.Range("D31").Value = 44490
$ws.Range("J31").Value = 72

$ws.Range("D32").Value = 44495
$ws.Range("J32").Value = 520
$ws.Range("K32").Value = 1800
$ws.Range("M32").Value = 1800
$ws.Range("P32").Value = 1800

$ws.Range("D33").Value = 44476
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 2000
$ws.Range("M33").Value = 2000
$ws.Range("P33").Value = 2000

$ws.Range("D34").Value = 44498
$ws.Range("J34").Value = 500
$ws.Range("K34").Value = 1600
$ws.Range("L34").Value = 1600
$ws.Range("M34").Value = 1600
$ws.Range("P34").Value = 1600

$ws.Range("D35").Value = 44159
$ws.Range("H35").Value = "Verde"
$ws.Range("J35").Value = 600
$ws.Range("L35").Value = 1700
$ws.Range("M35").Value = 1650
$ws.Range("P35").Value = 1650

$ws.Range("D36").Value = 44166
$ws.Range("H36").Value = "Verde"
$ws.Range("J36").Value = 500
$ws.Range("K36").Value = 1600
$ws.Range("L36").Value = 1600
$ws.Range("M36").Value = 1600
$ws.Range("P36").Value = 1600

$ws.Range("D37").Value = 44511
$ws.Range("J37").Value = 144
$ws.Range("K37").Value = 1700
$ws.Range("L37").Value = 1700

$ws.Range("D38").Value = 44481
$ws.Range("J38").Value = 300
$ws.Range("K38").Value = 1700
$ws.Range("L38").Value = 2000
$ws.Range("M38").Value = 1850
$ws.Range("P38").Value = 1850

$ws.Range("D39").Value = 44161
$ws.Range("H39").Value = "Verde"
$ws.Range("J39").Value = 300

$ws.Range("D40").Value = 44484
$ws.Range("J40").Value = 550
$ws.Range("K40").Value = 1700
$ws.Range("L40").Value = 1700
$ws.Range("M40").Value = 1700
$ws.Range("P40").Value = 1700

$ws.Range("D41").Value = 44526
$ws.Range("J41").Value = 600

$ws.Range("D42").Value = 44529
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("J42").Value = 120
$ws.Range("K42").Value = 1700
$ws.Range("M42").Value = 1700
$ws.Range("P42").Value = 1700

$ws.Range("D43").Value = 44522
$ws.Range("J43").Value = 180

$ws.Range("D44").Value = 44496
$ws.Range("J44").Value = 84
$ws.Range("K44").Value = 1800
$ws.Range("L44").Value = 1800
$ws.Range("M44").Value = 1800
$ws.Range("P44").Value = 1800

$ws.Range("D45").Value = 44162
$ws.Range("H45").Value = "Verde"
$ws.Range("J45").Value = 700
$ws.Range("K45").Value = 1600
$ws.Range("L45").Value = 1600
$ws.Range("M45").Value = 1600
$ws.Range("P45").Value = 1600

$ws.Range("D46").Value = 44517
$ws.Range("J46").Value = 120
